$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "73.330.44"

$ws.Range("D3").Value = "2.598.58"
$ws.Range("E3").Value = "  +6.49%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "182.96"
$ws.Range("E5").Value = "  +12.68%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "579.38"
$ws.Range("E6").Value = "  +3.40%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("E8").Value = "  +3.53%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.195"
$ws.Range("E9").Value = "  +15.62%  "

$ws.Range("D10").Value = "2.597.88"
$ws.Range("E10").Value = "  +6.61%  "

$ws.Range("E11").Value = "  -0.11%  "

$ws.Range("E12").Value = "  +7.70%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.73"
$ws.Range("E13").Value = "  +2.64%  "

$ws.Range("D14").Value = "73.239.50"
$ws.Range("E14").Value = "  +6.37%  "

$ws.Range("D15").Value = "3.076.32"
$ws.Range("E15").Value = "  +6.67%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000186"
$ws.Range("E16").Value = "  +4.81%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.80"
$ws.Range("E17").Value = "  +11.00%  "

$ws.Range("D18").Value = "2.599.22"
$ws.Range("E18").Value = "  +6.76%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.89"
$ws.Range("E19").Value = "  +28.68%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.76"
$ws.Range("E20").Value = "  +11.53%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "369.50"
$ws.Range("E21").Value = "  +8.70%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.21"
$ws.Range("E22").Value = "  +14.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.06"
$ws.Range("E23").Value = "  +5.63%  "

$ws.Range("E24").Value = "  -0.02%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.34"
$ws.Range("E25").Value = "  +3.35%  "

$ws.Range("E26").Value = "  +10.35%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.18"

$ws.Range("D28").Value = "2.713.34"
$ws.Range("E28").Value = "  +5.83%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.05%  "

$ws.Range("D30").Value = "0.0₃0923"
$ws.Range("E30").Value = "  +12.02%  "

$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "501.48"
$ws.Range("E31").Value = "  +16.82%  "

$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.36"
$ws.Range("E32").Value = "  +16.70%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.60"
$ws.Range("E33").Value = "  +6.41%  "

$ws.Range("E34").Value = "  +6.95%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.120"
$ws.Range("E36").Value = "  +12.67%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "160.31"
$ws.Range("E37").Value = "  +0.76%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.00"
$ws.Range("E38").Value = "  +5.73%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.22"
$ws.Range("E39").Value = "  +1.03%  "

$ws.Range("E40").Value = "  -0.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.82"
$ws.Range("E41").Value = "  +10.64%  "

$ws.Range("E42").Value = "  +9.36%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.321"
$ws.Range("E43").Value = "  +7.40%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "155.86"
$ws.Range("E44").Value = "  +20.16%  "

$ws.Range("E45").Value = "  +19.50%  "

$ws.Range("E46").Value = "  +8.34%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.32"
$ws.Range("E47").Value = "  +13.16%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "38.35"
$ws.Range("E48").Value = "  +2.46%  "

$ws.Range("E49").Value = "  +7.33%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.521"
$ws.Range("E50").Value = "  +8.26%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.08"
$ws.Range("E51").Value = "  +18.94%  "
